# Update the "names" sheet (Sheet1): remove the first id (row 1) by
# deleting the entire row, which shifts all remaining ids up by one.
$wb = $excel.ActiveWorkbook
$namesSheet = $wb.Worksheets.Item(1)
$usedSheet = $wb.Worksheets.Item(2)

# Capture the id that is about to be removed from the "names" sheet so it
# can be appended to the "used" sheet log below.
$removedId = $namesSheet.Range("A1").Value()

$namesSheet.Rows("1:1").Delete()

# Record the newly-used id in the "used" sheet, appending a new row right
# after the last existing entry.
$usedRange = $usedSheet.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$newRow = $lastRow + 1

$usedSheet.Cells.Item($newRow, 1).Value = $removedId
$usedSheet.Cells.Item($newRow, 2).Value = "ChatGPT Image 2026年1月18日 07_08_24.png"
$usedSheet.Cells.Item($newRow, 3).Value = "2026-01-18 07:09:47"
